# Updated cryptos list - applies latest price/volume figures and
# re-orders the MultiversX / FirstDigitalUSD rows (44 & 45).
#
# All of the touched cells hold plain text values (prices such as
# "42.724.56" use '.' as a thousands separator and are NOT valid
# numbers, and percentages are stored as literal strings padded with
# spaces, e.g. "  -1.05%  "). Assigning such strings straight to
# .Value works for values Excel cannot parse as a number, but plain
# numeric-looking strings like "1.00" or "68.19" would silently be
# converted into real numbers. To keep every one of these cells as
# text (matching the original workbook's inline-string cells) we
# briefly mark the cell as Text ("@") before writing the value and
# then clear the formatting again so no stray number format/style is
# left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $addr, $text)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
    $cell.ClearFormats()
}

Set-CellText $ws 'D2' '42.724.56'
Set-CellText $ws 'E2' '  -1.05%  '
Set-CellText $ws 'D3' '2.372.10'
Set-CellText $ws 'E3' '  +1.75%  '
Set-CellText $ws 'E4' '  -0.22%  '
Set-CellText $ws 'D5' '327.36'
Set-CellText $ws 'E5' '  +5.41%  '
Set-CellText $ws 'D6' '98.79'
Set-CellText $ws 'E6' '  -9.70%  '
Set-CellText $ws 'E7' '  +0.06%  '
Set-CellText $ws 'E8' '  +0.14%  '
Set-CellText $ws 'D9' '0.620'
Set-CellText $ws 'E9' '  +0.07%  '
Set-CellText $ws 'D10' '39.74'
Set-CellText $ws 'E10' '  -9.77%  '
Set-CellText $ws 'D11' '0.0920'
Set-CellText $ws 'E11' '  -1.28%  '
Set-CellText $ws 'D12' '8.38'
Set-CellText $ws 'E12' '  -5.79%  '
Set-CellText $ws 'E13' '  -6.21%  '
Set-CellText $ws 'E14' '  +0.30%  '
Set-CellText $ws 'D15' '16.28'
Set-CellText $ws 'E15' '  +2.48%  '
Set-CellText $ws 'D16' '2.731.94'
Set-CellText $ws 'E16' '  +2.65%  '
Set-CellText $ws 'D17' '2.379.20'
Set-CellText $ws 'E17' '  -1.22%  '
Set-CellText $ws 'D18' '42.726.60'
Set-CellText $ws 'E18' '  -0.99%  '
Set-CellText $ws 'D19' '7.81'
Set-CellText $ws 'E19' '  +6.87%  '
Set-CellText $ws 'E20' '  -2.23%  '
Set-CellText $ws 'D21' '3.71'
Set-CellText $ws 'E21' '  +6.46%  '
Set-CellText $ws 'D22' '75.06'
Set-CellText $ws 'E22' '  -1.24%  '
Set-CellText $ws 'D23' '271.23'
Set-CellText $ws 'E23' '  +6.28%  '
Set-CellText $ws 'D24' '2.31'
Set-CellText $ws 'E24' '  -8.33%  '
Set-CellText $ws 'D25' '9.67'
Set-CellText $ws 'E25' '  +6.47%  '
Set-CellText $ws 'E26' '  -0.07%  '
Set-CellText $ws 'D27' '11.41'
Set-CellText $ws 'E27' '  -4.10%  '
Set-CellText $ws 'D28' '23.58'
Set-CellText $ws 'E28' '  +4.46%  '
Set-CellText $ws 'E29' '  -1.69%  '
Set-CellText $ws 'D30' '172.65'
Set-CellText $ws 'E30' '  -0.79%  '
Set-CellText $ws 'D31' '3.11'
Set-CellText $ws 'E31' '  -1.74%  '
Set-CellText $ws 'D32' '0.0897'
Set-CellText $ws 'E32' '  -1.28%  '
Set-CellText $ws 'D33' '35.03'
Set-CellText $ws 'D34' '5.86'
Set-CellText $ws 'E34' '  +0.89%  '
Set-CellText $ws 'D35' '0.131'
Set-CellText $ws 'E35' '  -0.33%  '
Set-CellText $ws 'D36' '4.57'
Set-CellText $ws 'E36' '  -9.17%  '
Set-CellText $ws 'D37' '0.0356'
Set-CellText $ws 'E37' '  -5.88%  '
Set-CellText $ws 'D38' '3.86'
Set-CellText $ws 'E38' '  -7.34%  '
Set-CellText $ws 'D39' '0.104'
Set-CellText $ws 'E39' '  -0.31%  '
Set-CellText $ws 'E40' '  +3.40%  '
Set-CellText $ws 'E41' '  +0.94%  '
Set-CellText $ws 'D42' '95.90'
Set-CellText $ws 'E42' '  +50.37%  '
Set-CellText $ws 'E43' '  -3.54%  '
Set-CellText $ws 'B44' 'FirstDigitalUSD'
Set-CellText $ws 'C44' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-CellText $ws 'D44' '1.00'
Set-CellText $ws 'E44' '  +0.03%  '
Set-CellText $ws 'B45' 'MultiversX'
Set-CellText $ws 'C45' 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-CellText $ws 'D45' '68.19'
Set-CellText $ws 'E45' '  -5.17%  '
Set-CellText $ws 'D46' '115.84'
Set-CellText $ws 'E46' '  +4.63%  '
Set-CellText $ws 'D47' '11.74'
Set-CellText $ws 'E47' '  -6.04%  '
Set-CellText $ws 'D48' '5.40'
Set-CellText $ws 'E48' '  -5.98%  '
Set-CellText $ws 'D49' '8.92'
Set-CellText $ws 'E49' '  -2.50%  '
Set-CellText $ws 'D50' '1.603.14'
Set-CellText $ws 'E50' '  +7.94%  '
Set-CellText $ws 'E51' '  -3.03%  '
